# Daily Update 키워드 10개
# Strip literal "<b>" / "</b>" markup that leaked into product title strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B21" = "NIIMBOT D11 네임스티커 라벨프린터"
    "B29" = "HP 2131 가정용 프린터기 무한잉크 복합기 잉크젯 프린트 복사 공기계"
    "B30" = "아이티씨 퓨리웨이 휴대용 무선 라벨프린터 D11"
    "B34" = "디앤에스글로벌 휴대용 멀티 스티커 프린터 프릭커"
    "B35" = "캐논코리아정품 Canon SELPHY CP1300 셀프 사진 인화기 포토 프린터"
    "B38" = "엡손 EcoTank 흑백 프린터 M1120 정품 무한잉크"
    "B65" = "디엠케이인터내셔날 시노링크 라벨프린터 D11"
    "B71" = "엡손 정품 휴대용 라벨 프린터 네임스티커 라벨기 카카오 리락쿠마"
    "B75" = "삼성 프린터 컬러 레이저 프린트 가정용 사무실 무선 와이파이 레이져 복합기 팩스 스캔 복사"
    "B77" = "조이풀웍스 SIRO 휴대용 무선 라벨프린터 E10"
    "B79" = "HP8710 HP 9010 HP7740 HP7720 epson2850 가정용프린터 무한잉크 A3복합기"
    "B81" = "부라더 라벨프린터 모음[브라더 라벨기 휴대용 라벨터치 네임스티커 프린터]"
    "B84" = "삼성 복합기렌탈 CLX-6260FR 사무용 컬러 레이저 프린터 임대 36개월"
    "B87" = "hp오피스젯 7740 A3복합기 무한잉크 복합기 프린터 HP7740 HP7612 후속"
    "B89" = "샤오미 미지아 포켓 포토 프린터"
    "B91" = "캐논코리아정품 셀피 스퀘어 SELPHY SQUARE QX10 + 밤토리패키지 세트/ 휴대용 포토프린터 핸드폰 사진 인화기"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
